$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6804.4
$ws.Range("I43").Value = 12900
$ws.Range("J43").Value = 6127.1113
$ws.Range("K43").Value = 12900
$ws.Range("L43").Value = 6127.1113
$ws.Range("M43").Value = -12831
$ws.Range("N43").Value = -6265.1113

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1250
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 1000
$ws.Range("M86").Value = 123

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1250
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 5000
$ws.Range("M89").Value = 616

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7587.0967
$ws.Range("I116").Value = 9157.35
$ws.Range("J116").Value = 4732.091
$ws.Range("K116").Value = 9157.35
$ws.Range("L116").Value = 4732.091
$ws.Range("M116").Value = -5715.35
$ws.Range("N116").Value = -11616.091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1453.8043
$ws.Range("I132").Value = 1453.8043
$ws.Range("K132").Value = 4361.4129
$ws.Range("M132").Value = -1831.4129

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 93064.91
$ws.Range("I135").Value = 2371.5
$ws.Range("K135").Value = 21343.5
$ws.Range("M135").Value = -18808.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1934.2963
$ws.Range("I138").Value = 1493.1818
$ws.Range("J138").Value = 2237.5625
$ws.Range("K138").Value = 4479.5454
$ws.Range("L138").Value = 6712.6875
$ws.Range("M138").Value = 660.4546
$ws.Range("N138").Value = -16992.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5905.311
$ws.Range("I32").Value = 2480.2817
$ws.Range("K32").Value = 2480.2817
$ws.Range("M32").Value = -2193.2817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2057.4443
$ws.Range("J74").Value = 1987
$ws.Range("L74").Value = 1987
$ws.Range("N74").Value = -3735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2057.4443
$ws.Range("J77").Value = 1987
$ws.Range("L77").Value = 9935
$ws.Range("N77").Value = -18671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 112445.555
$ws.Range("I107").Value = 112445.555
$ws.Range("K107").Value = 112445.555
$ws.Range("M107").Value = -110525.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 27222
$ws.Range("J139").Value = 27222
$ws.Range("L139").Value = 27222
$ws.Range("N139").Value = -37502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 64111.8
$ws.Range("J141").Value = 64111.8
$ws.Range("L141").Value = 64111.8
$ws.Range("N141").Value = -74471.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32941.184
$ws.Range("I31").Value = 31163.457
$ws.Range("J31").Value = 39854.555
$ws.Range("K31").Value = 31163.457
$ws.Range("L31").Value = 39854.555
$ws.Range("M31").Value = -30868.457
$ws.Range("N31").Value = -40444.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 32941.184
$ws.Range("I34").Value = 31163.457
$ws.Range("J34").Value = 39854.555
$ws.Range("K34").Value = 31163.457
$ws.Range("L34").Value = 39854.555
$ws.Range("M34").Value = -30961.457
$ws.Range("N34").Value = -40258.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 45935
$ws.Range("J52").Value = 49996.5
$ws.Range("L52").Value = 49996.5
$ws.Range("N52").Value = -50584.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3576.6086
$ws.Range("I132").Value = 3446.0952
$ws.Range("K132").Value = 10338.2856
$ws.Range("M132").Value = -7808.285600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10240
$ws.Range("I134").Value = 10600
$ws.Range("K134").Value = 31800
$ws.Range("M134").Value = -29265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 152780
$ws.Range("J135").Value = 152780
$ws.Range("L135").Value = 152780
$ws.Range("N135").Value = -162920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 238992.36
$ws.Range("J141").Value = 238992.36
$ws.Range("L141").Value = 238992.36
$ws.Range("N141").Value = -249352.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2647359
$ws.Range("I4").Value = 1111611.2
$ws.Range("J4").Value = 6102791.5
$ws.Range("K4").Value = 3334833.6
$ws.Range("L4").Value = 18308374.5
$ws.Range("M4").Value = -3334721.6
$ws.Range("N4").Value = -18308598.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1899.5
$ws.Range("J64").Value = 1899.5
$ws.Range("L64").Value = 5698.5
$ws.Range("N64").Value = -6238.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 1899.5
$ws.Range("J67").Value = 1899.5
$ws.Range("L67").Value = 5698.5
$ws.Range("N67").Value = -7570.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2129.4119
$ws.Range("J131").Value = 2011.3334
$ws.Range("L131").Value = 6034.0002
$ws.Range("N131").Value = -16114.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3011.6667
$ws.Range("I137").Value = 1887.3334
$ws.Range("K137").Value = 5662.0002
$ws.Range("M137").Value = -562.0002000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2184.7
$ws.Range("I138").Value = 1979.3334
$ws.Range("K138").Value = 5938.0002
$ws.Range("M138").Value = -798.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1270.5
$ws.Range("J9").Value = 2000
$ws.Range("L9").Value = 2000
$ws.Range("N9").Value = -2340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 225687.75
$ws.Range("I16").Value = 160900
$ws.Range("J16").Value = 333667.34
$ws.Range("K16").Value = 160900
$ws.Range("L16").Value = 333667.34
$ws.Range("M16").Value = -160730
$ws.Range("N16").Value = -334007.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5305.75
$ws.Range("I136").Value = 4240.5713
$ws.Range("K136").Value = 12721.7139
$ws.Range("M136").Value = -10171.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1773.3077
$ws.Range("I81").Value = 1773.3077
$ws.Range("K81").Value = 3546.6154
$ws.Range("M81").Value = -2485.6154

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1773.3077
$ws.Range("I84").Value = 1773.3077
$ws.Range("K84").Value = 17733.077
$ws.Range("M84").Value = -12429.077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3162.0476
$ws.Range("I132").Value = 3391.75
$ws.Range("K132").Value = 10175.25
$ws.Range("M132").Value = -7645.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 140000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 140000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 140000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -150280
